# "agrego cambios de pedidos"
# - Add a new "Region Metropolitana (sede CABA) -> Buenos Aires" line under the
#   existing region descriptions (B18:B20 -> now B18:B21).
# - Clarify the Local/Provincial/Regional/Nacional column headers in the first
#   pricing table (row 3) with short parenthetical explanations.
# - Re-fit columns C:F so the new, longer header text is fully visible.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New region row, added right after "Region Norte..." (B20) -> becomes B21.
$ws.Range("B21").Value = "Region Metropolitana (sede CABA) -> Buenos Aires"

# Clarify the pricing-table column headers (row 3, columns C:F).
$ws.Range("C3").Value = "Local (misma localidad)"
$ws.Range("D3").Value = "Provincial (misma provincia)"
$ws.Range("E3").Value = "Regional (misma región)"
$ws.Range("F3").Value = "Nacional (inter-regional)"

# Widen columns C:F to fit the new, longer header text.
$ws.Columns("C:F").AutoFit() | Out-Null

# Match the new selection left behind in the source file.
$ws.Range("C18:C21").Select() | Out-Null
